# Efnb2-Grm5.xlsx -- update scripts with new TPM values
# Re-computed NATMI LR-pair edge weights (rows 2-11) and added three
# new Sending-cluster x Target-cluster rows (12-16) for the MuSCs /
# Resolving-Mac x Inflammatory-Mac / Resolving-Mac combinations.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (updated row)
$ws.Range("D2").Value = "FAPs"
$ws.Range("G2").Value = 37.98277566666666
$ws.Range("H2").Value = 113.948327
$ws.Range("I2").Value = 0.697850645410475
$ws.Range("J2").Value = 0.6978506454104751
$ws.Range("M2").Value = 0.071834
$ws.Range("N2").Value = 0.215502
$ws.Range("O2").Value = 0.9835064532028698
$ws.Range("P2").Value = 0.9835064532028697
$ws.Range("Q2").Value = 2.728454707239333
$ws.Range("R2").Value = 24.556092365154
$ws.Range("S2").Value = 0.6863406131329899
$ws.Range("T2").Value = 0.6863406131329899

# Row 3 (updated row)
$ws.Range("D3").Value = "Inflammatory-Mac"
$ws.Range("G3").Value = 37.98277566666666
$ws.Range("H3").Value = 113.948327
$ws.Range("I3").Value = 0.697850645410475
$ws.Range("J3").Value = 0.6978506454104751
$ws.Range("M3").Value = 0.0002013333333333333
$ws.Range("N3").Value = 0.000604
$ws.Range("O3").Value = 0.002756530787345516
$ws.Range("P3").Value = 0.002756530787345516
$ws.Range("Q3").Value = 0.007647198834222222
$ws.Range("R3").Value = 0.068824789508
$ws.Range("S3").Value = 0.001923646789042913
$ws.Range("T3").Value = 0.001923646789042913

# Row 4 (updated row)
$ws.Range("A4").Value = "ECs"
$ws.Range("D4").Value = "Resolving-Mac"
$ws.Range("G4").Value = 37.98277566666666
$ws.Range("H4").Value = 113.948327
$ws.Range("I4").Value = 0.697850645410475
$ws.Range("J4").Value = 0.6978506454104751
$ws.Range("M4").Value = 0.001003333333333333
$ws.Range("N4").Value = 0.00301
$ws.Range("O4").Value = 0.01373701600978477
$ws.Range("P4").Value = 0.01373701600978477
$ws.Range("Q4").Value = 0.03810938491888888
$ws.Range("R4").Value = 0.34298446427
$ws.Range("S4").Value = 0.009586385488442333
$ws.Range("T4").Value = 0.009586385488442333

# Row 5 (updated row)
$ws.Range("I5").Value = 0.1779541659542351
$ws.Range("J5").Value = 0.1779541659542352
$ws.Range("O5").Value = 0.9835064532028698
$ws.Range("P5").Value = 0.9835064532028697
$ws.Range("S5").Value = 0.1750190705903247
$ws.Range("T5").Value = 0.1750190705903247

# Row 6 (updated row)
$ws.Range("A6").Value = "FAPs"
$ws.Range("D6").Value = "Inflammatory-Mac"
$ws.Range("G6").Value = 9.685730333333334
$ws.Range("H6").Value = 29.057191
$ws.Range("I6").Value = 0.1779541659542351
$ws.Range("J6").Value = 0.1779541659542352
$ws.Range("M6").Value = 0.0002013333333333333
$ws.Range("N6").Value = 0.000604
$ws.Range("O6").Value = 0.002756530787345516
$ws.Range("P6").Value = 0.002756530787345516
$ws.Range("Q6").Value = 0.001950060373777778
$ws.Range("R6").Value = 0.017550543364
$ws.Range("S6").Value = 0.0004905361371892425
$ws.Range("T6").Value = 0.0004905361371892425

# Row 7 (updated row)
$ws.Range("A7").Value = "FAPs"
$ws.Range("D7").Value = "Resolving-Mac"
$ws.Range("G7").Value = 9.685730333333334
$ws.Range("H7").Value = 29.057191
$ws.Range("I7").Value = 0.1779541659542351
$ws.Range("J7").Value = 0.1779541659542352
$ws.Range("M7").Value = 0.001003333333333333
$ws.Range("N7").Value = 0.00301
$ws.Range("O7").Value = 0.01373701600978477
$ws.Range("P7").Value = 0.01373701600978477
$ws.Range("Q7").Value = 0.009718016101111112
$ws.Range("R7").Value = 0.08746214491000001
$ws.Range("S7").Value = 0.002444559226721224
$ws.Range("T7").Value = 0.002444559226721224

# Row 8 (updated row)
$ws.Range("A8").Value = "Inflammatory-Mac"
$ws.Range("D8").Value = "FAPs"
$ws.Range("G8").Value = 0.5676613333333332
$ws.Range("H8").Value = 1.702984
$ws.Range("I8").Value = 0.01042953867610283
$ws.Range("J8").Value = 0.01042953867610283
$ws.Range("M8").Value = 0.071834
$ws.Range("N8").Value = 0.215502
$ws.Range("O8").Value = 0.9835064532028698
$ws.Range("P8").Value = 0.9835064532028697
$ws.Range("Q8").Value = 0.04077738421866665
$ws.Range("R8").Value = 0.366996457968
$ws.Range("S8").Value = 0.01025751859187605
$ws.Range("T8").Value = 0.01025751859187605

# Row 9 (updated row)
$ws.Range("A9").Value = "Inflammatory-Mac"
$ws.Range("D9").Value = "Inflammatory-Mac"
$ws.Range("G9").Value = 0.5676613333333332
$ws.Range("H9").Value = 1.702984
$ws.Range("I9").Value = 0.01042953867610283
$ws.Range("J9").Value = 0.01042953867610283
$ws.Range("M9").Value = 0.0002013333333333333
$ws.Range("N9").Value = 0.000604
$ws.Range("O9").Value = 0.002756530787345516
$ws.Range("P9").Value = 0.002756530787345516
$ws.Range("Q9").Value = 0.0001142891484444444
$ws.Range("R9").Value = 0.001028602336
$ws.Range("S9").Value = 0.00002874934445848825
$ws.Range("T9").Value = 0.00002874934445848825

# Row 10 (updated row)
$ws.Range("A10").Value = "Inflammatory-Mac"
$ws.Range("D10").Value = "Resolving-Mac"
$ws.Range("G10").Value = 0.5676613333333332
$ws.Range("H10").Value = 1.702984
$ws.Range("I10").Value = 0.01042953867610283
$ws.Range("J10").Value = 0.01042953867610283
$ws.Range("M10").Value = 0.001003333333333333
$ws.Range("N10").Value = 0.00301
$ws.Range("O10").Value = 0.01373701600978477
$ws.Range("P10").Value = 0.01373701600978477
$ws.Range("Q10").Value = 0.0005695535377777777
$ws.Range("R10").Value = 0.00512598184
$ws.Range("S10").Value = 0.0001432707397682941
$ws.Range("T10").Value = 0.0001432707397682941

# Row 11 (updated row)
$ws.Range("A11").Value = "MuSCs"
$ws.Range("G11").Value = 5.823095333333334
$ws.Range("H11").Value = 17.469286
$ws.Range("I11").Value = 0.1069866739681064
$ws.Range("J11").Value = 0.1069866739681064
$ws.Range("O11").Value = 0.9835064532028698
$ws.Range("P11").Value = 0.9835064532028697
$ws.Range("Q11").Value = 0.4182962301746667
$ws.Range("R11").Value = 3.764666071572
$ws.Range("S11").Value = 0.1052220842543442
$ws.Range("T11").Value = 0.1052220842543442

# Row 12 (new row)
$ws.Range("A12").Value = "MuSCs"
$ws.Range("B12").Value = "Efnb2"
$ws.Range("C12").Value = "Grm5"
$ws.Range("D12").Value = "Inflammatory-Mac"
$ws.Range("E12").Value = 3
$ws.Range("F12").Value = 1
$ws.Range("G12").Value = 5.823095333333334
$ws.Range("H12").Value = 17.469286
$ws.Range("I12").Value = 0.1069866739681064
$ws.Range("J12").Value = 0.1069866739681064
$ws.Range("K12").Value = 1
$ws.Range("L12").Value = 0.3333333333333333
$ws.Range("M12").Value = 0.0002013333333333333
$ws.Range("N12").Value = 0.000604
$ws.Range("O12").Value = 0.002756530787345516
$ws.Range("P12").Value = 0.002756530787345516
$ws.Range("Q12").Value = 0.001172383193777778
$ws.Range("R12").Value = 0.010551448744
$ws.Range("S12").Value = 0.0002949120606287825
$ws.Range("T12").Value = 0.0002949120606287824

# Row 13 (new row)
$ws.Range("A13").Value = "MuSCs"
$ws.Range("B13").Value = "Efnb2"
$ws.Range("C13").Value = "Grm5"
$ws.Range("D13").Value = "Resolving-Mac"
$ws.Range("E13").Value = 3
$ws.Range("F13").Value = 1
$ws.Range("G13").Value = 5.823095333333334
$ws.Range("H13").Value = 17.469286
$ws.Range("I13").Value = 0.1069866739681064
$ws.Range("J13").Value = 0.1069866739681064
$ws.Range("K13").Value = 1
$ws.Range("L13").Value = 0.3333333333333333
$ws.Range("M13").Value = 0.001003333333333333
$ws.Range("N13").Value = 0.00301
$ws.Range("O13").Value = 0.01373701600978477
$ws.Range("P13").Value = 0.01373701600978477
$ws.Range("Q13").Value = 0.005842505651111112
$ws.Range("R13").Value = 0.05258255086
$ws.Range("S13").Value = 0.001469677653133502
$ws.Range("T13").Value = 0.001469677653133502

# Row 14 (new row)
$ws.Range("A14").Value = "Resolving-Mac"
$ws.Range("B14").Value = "Efnb2"
$ws.Range("C14").Value = "Grm5"
$ws.Range("D14").Value = "FAPs"
$ws.Range("E14").Value = 3
$ws.Range("F14").Value = 1
$ws.Range("G14").Value = 0.3689676666666666
$ws.Range("H14").Value = 1.106903
$ws.Range("I14").Value = 0.006778975991080511
$ws.Range("J14").Value = 0.006778975991080512
$ws.Range("K14").Value = 1
$ws.Range("L14").Value = 0.3333333333333333
$ws.Range("M14").Value = 0.071834
$ws.Range("N14").Value = 0.215502
$ws.Range("O14").Value = 0.9835064532028698
$ws.Range("P14").Value = 0.9835064532028697
$ws.Range("Q14").Value = 0.02650442336733333
$ws.Range("R14").Value = 0.238539810306
$ws.Range("S14").Value = 0.006667166633335002
$ws.Range("T14").Value = 0.006667166633335002

# Row 15 (new row)
$ws.Range("A15").Value = "Resolving-Mac"
$ws.Range("B15").Value = "Efnb2"
$ws.Range("C15").Value = "Grm5"
$ws.Range("D15").Value = "Inflammatory-Mac"
$ws.Range("E15").Value = 3
$ws.Range("F15").Value = 1
$ws.Range("G15").Value = 0.3689676666666666
$ws.Range("H15").Value = 1.106903
$ws.Range("I15").Value = 0.006778975991080511
$ws.Range("J15").Value = 0.006778975991080512
$ws.Range("K15").Value = 1
$ws.Range("L15").Value = 0.3333333333333333
$ws.Range("M15").Value = 0.0002013333333333333
$ws.Range("N15").Value = 0.000604
$ws.Range("O15").Value = 0.002756530787345516
$ws.Range("P15").Value = 0.002756530787345516
$ws.Range("Q15").Value = 0.00007428549022222222
$ws.Range("R15").Value = 0.000668569412
$ws.Range("S15").Value = 0.00001868645602608951
$ws.Range("T15").Value = 0.00001868645602608951

# Row 16 (new row)
$ws.Range("A16").Value = "Resolving-Mac"
$ws.Range("B16").Value = "Efnb2"
$ws.Range("C16").Value = "Grm5"
$ws.Range("D16").Value = "Resolving-Mac"
$ws.Range("E16").Value = 3
$ws.Range("F16").Value = 1
$ws.Range("G16").Value = 0.3689676666666666
$ws.Range("H16").Value = 1.106903
$ws.Range("I16").Value = 0.006778975991080511
$ws.Range("J16").Value = 0.006778975991080512
$ws.Range("K16").Value = 1
$ws.Range("L16").Value = 0.3333333333333333
$ws.Range("M16").Value = 0.001003333333333333
$ws.Range("N16").Value = 0.00301
$ws.Range("O16").Value = 0.01373701600978477
$ws.Range("P16").Value = 0.01373701600978477
$ws.Range("Q16").Value = 0.0003701975588888888
$ws.Range("R16").Value = 0.00333177803
$ws.Range("S16").Value = 0.00009312290171941957
$ws.Range("T16").Value = 0.00009312290171941957
